$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "maxiter" row (row 8) values from 2000 to 20 across columns B:S
$ws.Range("B8:S8").Value = 20

# Update the active cell selection to F12
$ws.Range("F12").Select()
